# Update slide 13 ("Seattle Housing Prices - Demonstration"):
#  - Split the body placeholder's first line so the word "tkinter" is its
#    own run, and append several new paragraphs describing the newly
#    imported libraries (with a couple of bold, second-level code lines).
#  - Remove the separate red-outlined "Update" flag textbox that used to
#    sit in the top-right corner of the slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)

$body = $s.Shapes.Item(2)
$tr = $body.TextFrame.TextRange

$line1 = "Visual demonstration of how to estimate housing prices using the selected metrics using the tkinter graphical user interface   "
$line2 = "Imported new libraries "
$line3 = "import tkinter as tk"
$line4 = "from tkinter import scrolledtext"
$line5 = "from tkinter.scrolledtext import"
$line6 = ""
$line7 = " "

$tr.Text = $line1 + "`r" + $line2 + "`r" + $line3 + "`r" + $line4 + "`r" + $line5 + "`r" + $line6 + "`r" + $line7

# Helper: re-stamp a sub-range's own font size onto itself so the run
# boundary is preserved (needed because PowerPoint marks "tkinter" /
# "scrolledtext" / "tk" as flagged-misspelling runs even though their
# visible formatting matches the text around them).
function Split-Run($range) {
    $range.Font.Size = $range.Font.Size
}

# --- Paragraph 1: "...using the tkinter graphical user interface" ---
# Carve the word "tkinter" out into its own run (no other formatting change;
# it is just flagged by the spell-checker as an unrecognized word).
$para1 = $tr.Paragraphs(1)
$tkinterStart1 = $line1.IndexOf("tkinter") + 1
Split-Run $para1.Characters($tkinterStart1, 7)

# --- Paragraph 2: "Imported new libraries " (no extra formatting) ---

# --- Paragraph 3: "import tkinter as tk" -> bold, second outline level ---
$para3 = $tr.Paragraphs(3)
$para3.IndentLevel = 2
$para3.Font.Bold = $true
Split-Run $para3.Characters(($line3.IndexOf("tkinter") + 1), 7)
Split-Run $para3.Characters(($line3.LastIndexOf("tk") + 1), 2)

# --- Paragraph 4: "from tkinter import scrolledtext" -> bold, level 2 ---
$para4 = $tr.Paragraphs(4)
$para4.IndentLevel = 2
$para4.Font.Bold = $true
Split-Run $para4.Characters(($line4.IndexOf("tkinter") + 1), 7)
Split-Run $para4.Characters(($line4.IndexOf("scrolledtext") + 1), 12)

# --- Paragraph 5: "from tkinter.scrolledtext import" -> bold, level 2 ---
$para5 = $tr.Paragraphs(5)
$para5.IndentLevel = 2
$para5.Font.Bold = $true
Split-Run $para5.Characters(($line5.IndexOf("tkinter.scrolledtext") + 1), 20)

# --- Paragraph 6: blank spacer line (plain formatting) ---

# --- Paragraph 7: trailing blank line with no bullet / no indent ---
$para7 = $tr.Paragraphs(7)
$para7.ParagraphFormat.Bullet.Visible = $false

# Remove the red "Update" flag textbox (was Shapes.Item(3) / id=5).
$s.Shapes.Item(3).Delete()
